# Auto-generated edit script: refresh market-price derived values
# in the FFXIV Leve profit workbook (scheduled runner data update).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H80").Value = 2880.7
$ws.Range("I80").Value = 1174.2307
$ws.Range("K80").Value = 3522.6921
$ws.Range("M80").Value = -2524.6921
$ws.Range("H83").Value = 2880.7
$ws.Range("I83").Value = 1174.2307
$ws.Range("K83").Value = 10568.0763
$ws.Range("M83").Value = -5576.076300000001
$ws.Range("H106").Value = 2116.5386
$ws.Range("I106").Value = 1866.5
$ws.Range("J106").Value = 2950
$ws.Range("K106").Value = 1866.5
$ws.Range("L106").Value = 2950
$ws.Range("M106").Value = -1235.5
$ws.Range("N106").Value = -4212
$ws.Range("H112").Value = 1349.2963
$ws.Range("J112").Value = 1333.375
$ws.Range("L112").Value = 4000.125
$ws.Range("N112").Value = -6216.125
$ws.Range("H130").Value = 65000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040
$ws.Range("H138").Value = 7911.7896
$ws.Range("I138").Value = 3048.75
$ws.Range("J138").Value = 8483.912
$ws.Range("K138").Value = 9146.25
$ws.Range("L138").Value = 25451.736
$ws.Range("M138").Value = -4006.25
$ws.Range("N138").Value = -35731.736

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 13894943
$ws.Range("I32").Value = 15391026
$ws.Range("K32").Value = 15391026
$ws.Range("M32").Value = -15390739
$ws.Range("H61").Value = 3160
$ws.Range("I61").Value = 2477.6843
$ws.Range("J61").Value = 5752.8
$ws.Range("K61").Value = 2477.6843
$ws.Range("L61").Value = 5752.8
$ws.Range("M61").Value = -2265.6843
$ws.Range("N61").Value = -6176.8
$ws.Range("H136").Value = 3160
$ws.Range("I136").Value = 2477.6843
$ws.Range("J136").Value = 5752.8
$ws.Range("K136").Value = 7433.0529
$ws.Range("L136").Value = 17258.4
$ws.Range("M136").Value = -4883.0529
$ws.Range("N136").Value = -22358.4
$ws.Range("H137").Value = 48749.5
$ws.Range("I137").Value = 47499
$ws.Range("K137").Value = 47499
$ws.Range("M137").Value = -42399

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 52582.7
$ws.Range("I86").Value = 60428.06
$ws.Range("J86").Value = 8125.6665
$ws.Range("K86").Value = 60428.06
$ws.Range("L86").Value = 8125.6665
$ws.Range("M86").Value = -59305.06
$ws.Range("N86").Value = -10371.6665
$ws.Range("H89").Value = 52582.7
$ws.Range("I89").Value = 60428.06
$ws.Range("J89").Value = 8125.6665
$ws.Range("K89").Value = 302140.3
$ws.Range("L89").Value = 40628.3325
$ws.Range("M89").Value = -296524.3
$ws.Range("N89").Value = -51860.3325
$ws.Range("H94").Value = 4276.3
$ws.Range("I94").Value = 3765.2666
$ws.Range("J94").Value = 5809.4
$ws.Range("K94").Value = 3765.2666
$ws.Range("L94").Value = 5809.4
$ws.Range("M94").Value = -3314.2666
$ws.Range("N94").Value = -6711.4

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 1250000
$ws.Range("I4").Value = 1250000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1250000
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("M4").Value = -1249888
$ws.Range("H31").Value = 2652.75
$ws.Range("I31").Value = 1975.9231
$ws.Range("J31").Value = 3909.7144
$ws.Range("K31").Value = 1975.9231
$ws.Range("L31").Value = 3909.7144
$ws.Range("M31").Value = -1680.9231
$ws.Range("N31").Value = -4499.7144
$ws.Range("H34").Value = 2652.75
$ws.Range("I34").Value = 1975.9231
$ws.Range("J34").Value = 3909.7144
$ws.Range("K34").Value = 1975.9231
$ws.Range("L34").Value = 3909.7144
$ws.Range("M34").Value = -1773.9231
$ws.Range("N34").Value = -4313.7144

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 1727.3334
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 9000
$ws.Range("N68").Value = -10622
$ws.Range("H71").Value = 1727.3334
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 27000
$ws.Range("N71").Value = -35112
$ws.Range("H88").Value = 4715.643
$ws.Range("I88").Value = 4125
$ws.Range("J88").Value = 4951.9
$ws.Range("K88").Value = 12375
$ws.Range("L88").Value = 14855.7
$ws.Range("M88").Value = -11947
$ws.Range("N88").Value = -15711.7
$ws.Range("H91").Value = 4715.643
$ws.Range("I91").Value = 4125
$ws.Range("J91").Value = 4951.9
$ws.Range("K91").Value = 12375
$ws.Range("L91").Value = 14855.7
$ws.Range("M91").Value = -10893
$ws.Range("N91").Value = -17819.7

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 5194.1875
$ws.Range("I70").Value = 4809.8335
$ws.Range("J70").Value = 6347.25
$ws.Range("K70").Value = 4809.8335
$ws.Range("L70").Value = 6347.25
$ws.Range("M70").Value = -4539.8335
$ws.Range("N70").Value = -6887.25
$ws.Range("H73").Value = 5194.1875
$ws.Range("I73").Value = 4809.8335
$ws.Range("J73").Value = 6347.25
$ws.Range("K73").Value = 4809.8335
$ws.Range("L73").Value = 6347.25
$ws.Range("M73").Value = -3873.8335
$ws.Range("N73").Value = -8219.25

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 3280.818
$ws.Range("I16").Value = 3308.7
$ws.Range("J16").Value = 3002
$ws.Range("K16").Value = 3308.7
$ws.Range("L16").Value = 3002
$ws.Range("N16").Value = -3342
$ws.Range("M16").Value = -3138.7
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = $null
$ws.Range("H46").Value = 2183.5293
$ws.Range("I46").Value = 977.1818
$ws.Range("K46").Value = 977.1818
$ws.Range("M46").Value = -789.1818
$ws.Range("H55").Value = 541.1429000000001
$ws.Range("I55").Value = 477.73334
$ws.Range("J55").Value = 699.6667
$ws.Range("K55").Value = 477.73334
$ws.Range("L55").Value = 699.6667
$ws.Range("M55").Value = -304.73334
$ws.Range("N55").Value = -1045.6667
$ws.Range("H132").Value = 5317.533
$ws.Range("I132").Value = 5026.0938
$ws.Range("J132").Value = 6034.923
$ws.Range("K132").Value = 15078.2814
$ws.Range("L132").Value = 18104.769
$ws.Range("M132").Value = -12548.2814
$ws.Range("N132").Value = -23164.769
$ws.Range("H136").Value = 4343.773
$ws.Range("I136").Value = 2681.5454
$ws.Range("J136").Value = 6006
$ws.Range("K136").Value = 8044.6362
$ws.Range("L136").Value = 18018
$ws.Range("M136").Value = -5494.6362
$ws.Range("N136").Value = -23118

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H70").Value = 24177.375
$ws.Range("I70").Value = 12500
$ws.Range("J70").Value = 28069.834
$ws.Range("K70").Value = 12500
$ws.Range("L70").Value = 28069.834
$ws.Range("M70").Value = -12185
$ws.Range("N70").Value = -28699.834
$ws.Range("H73").Value = 24177.375
$ws.Range("I73").Value = 12500
$ws.Range("J73").Value = 28069.834
$ws.Range("K73").Value = 12500
$ws.Range("L73").Value = 28069.834
$ws.Range("M73").Value = -11408
$ws.Range("N73").Value = -30253.834
$ws.Range("H100").Value = 432.9
$ws.Range("I100").Value = 400.47058
$ws.Range("J100").Value = 616.6667
$ws.Range("K100").Value = 800.94116
$ws.Range("L100").Value = 1233.3334
$ws.Range("M100").Value = -259.94116
$ws.Range("N100").Value = -2315.3334
$ws.Range("H123").Value = 32750
$ws.Range("J123").Value = 32750
$ws.Range("L123").Value = 32750
$ws.Range("N123").Value = -42550
$ws.Range("H138").Value = 59000
$ws.Range("J138").Value = 59000
$ws.Range("L138").Value = 59000
$ws.Range("N138").Value = -69280

